$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Completion status" column header in C2
$ws.Range("C2").Value = "完成情况"

# Mark rows 3-6 in column C as "已完成" (Completed)
$ws.Range("C3").Value = "已完成"
$ws.Range("C4").Value = "已完成"
$ws.Range("C5").Value = "已完成"
$ws.Range("C6").Value = "已完成"

# Update the active selection to G16
$ws.Range("G16").Select()
